# Auto-generated edit script: refreshes market-data columns (H:N)
# on sheets ARM, BSM, CRP, CUL, GSM, LTW, WVR to match the scheduled
# market-data runner's latest pull.

$wb = $excel.ActiveWorkbook

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5579.304
$ws.Range("I74").Value = 948
$ws.Range("K74").Value = 948
$ws.Range("M74").Value = -74
$ws.Range("H77").Value = 5579.304
$ws.Range("I77").Value = 948
$ws.Range("K77").Value = 4740
$ws.Range("M77").Value = -372
$ws.Range("H80").Value = 17492.8
$ws.Range("J80").Value = 17313.715
$ws.Range("L80").Value = 17313.715
$ws.Range("N80").Value = -19309.715
$ws.Range("H81").Value = 40000
$ws.Range("J81").Value = 40000
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -41996
$ws.Range("H82").Value = 39790
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 39790
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 39790
$ws.Range("M82").ClearContents() | Out-Null
$ws.Range("N82").Value = -40512
$ws.Range("H83").Value = 17492.8
$ws.Range("J83").Value = 17313.715
$ws.Range("L83").Value = 51941.145
$ws.Range("N83").Value = -61925.145
$ws.Range("H84").Value = 40000
$ws.Range("J84").Value = 40000
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -129984
$ws.Range("H85").Value = 39790
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 39790
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 39790
$ws.Range("M85").ClearContents() | Out-Null
$ws.Range("N85").Value = -42286
$ws.Range("H122").Value = 8334718
$ws.Range("I122").Value = 10527643
$ws.Range("J122").Value = 1600.2
$ws.Range("K122").Value = 31582929
$ws.Range("L122").Value = 4800.6
$ws.Range("M122").Value = -31580479
$ws.Range("N122").Value = -9700.6

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 13922.05
$ws.Range("I82").Value = 5778.5
$ws.Range("J82").Value = 15957.9375
$ws.Range("K82").Value = 5778.5
$ws.Range("L82").Value = 15957.9375
$ws.Range("M82").Value = -5395.5
$ws.Range("N82").Value = -16723.9375
$ws.Range("H85").Value = 13922.05
$ws.Range("I85").Value = 5778.5
$ws.Range("J85").Value = 15957.9375
$ws.Range("K85").Value = 5778.5
$ws.Range("L85").Value = 15957.9375
$ws.Range("M85").Value = -4452.5
$ws.Range("N85").Value = -18609.9375
$ws.Range("H134").Value = 20835214
$ws.Range("I134").Value = 25642556
$ws.Range("J134").Value = 3407.3333
$ws.Range("K134").Value = 76927668
$ws.Range("L134").Value = 10221.9999
$ws.Range("M134").Value = -76925133
$ws.Range("N134").Value = -15291.9999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13977.186
$ws.Range("I31").Value = 1342.8572
$ws.Range("K31").Value = 1342.8572
$ws.Range("M31").Value = -1047.8572
$ws.Range("H34").Value = 13977.186
$ws.Range("I34").Value = 1342.8572
$ws.Range("K34").Value = 1342.8572
$ws.Range("M34").Value = -1140.8572

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 627.5
$ws.Range("I2").Value = 26.285715
$ws.Range("J2").Value = 1010.0909
$ws.Range("K2").Value = 157.71429
$ws.Range("L2").Value = 6060.5454
$ws.Range("M2").Value = -44.71429000000001
$ws.Range("N2").Value = -6286.5454
$ws.Range("H41").Value = 725.9
$ws.Range("J41").Value = 599.6667
$ws.Range("L41").Value = 1799.0001
$ws.Range("N41").Value = -2475.0001
$ws.Range("H64").Value = 2666.3076
$ws.Range("I64").Value = 765.5
$ws.Range("J64").Value = 3511.111
$ws.Range("K64").Value = 2296.5
$ws.Range("L64").Value = 10533.333
$ws.Range("M64").Value = -2026.5
$ws.Range("N64").Value = -11073.333
$ws.Range("H67").Value = 2666.3076
$ws.Range("I67").Value = 765.5
$ws.Range("J67").Value = 3511.111
$ws.Range("K67").Value = 2296.5
$ws.Range("L67").Value = 10533.333
$ws.Range("M67").Value = -1360.5
$ws.Range("N67").Value = -12405.333
$ws.Range("H70").Value = 1459.1428
$ws.Range("I70").Value = 882.8
$ws.Range("J70").Value = 2900
$ws.Range("K70").Value = 2648.4
$ws.Range("L70").Value = 8700
$ws.Range("M70").Value = -2333.4
$ws.Range("N70").Value = -9330
$ws.Range("H73").Value = 1459.1428
$ws.Range("I73").Value = 882.8
$ws.Range("J73").Value = 2900
$ws.Range("K73").Value = 2648.4
$ws.Range("L73").Value = 8700
$ws.Range("M73").Value = -1556.4
$ws.Range("N73").Value = -10884
$ws.Range("H76").Value = 1998.3334
$ws.Range("I76").Value = 1747.5
$ws.Range("K76").Value = 5242.5
$ws.Range("M76").Value = -4859.5
$ws.Range("H79").Value = 1998.3334
$ws.Range("I79").Value = 1747.5
$ws.Range("K79").Value = 5242.5
$ws.Range("M79").Value = -3916.5
$ws.Range("H82").Value = 4357.143
$ws.Range("I82").Value = 2000
$ws.Range("J82").Value = 4750
$ws.Range("K82").Value = 6000
$ws.Range("L82").Value = 14250
$ws.Range("M82").Value = -5594
$ws.Range("N82").Value = -15062
$ws.Range("H85").Value = 4357.143
$ws.Range("I85").Value = 2000
$ws.Range("J85").Value = 4750
$ws.Range("K85").Value = 6000
$ws.Range("L85").Value = 14250
$ws.Range("M85").Value = -4596
$ws.Range("N85").Value = -17058
$ws.Range("H88").Value = 3660
$ws.Range("J88").Value = 3660
$ws.Range("L88").Value = 10980
$ws.Range("N88").Value = -11836
$ws.Range("H91").Value = 3660
$ws.Range("J91").Value = 3660
$ws.Range("L91").Value = 10980
$ws.Range("N91").Value = -13944

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 21585.492
$ws.Range("I70").Value = 31604.666
$ws.Range("J70").Value = 4409.7617
$ws.Range("K70").Value = 31604.666
$ws.Range("L70").Value = 4409.7617
$ws.Range("M70").Value = -31334.666
$ws.Range("N70").Value = -4949.7617
$ws.Range("H73").Value = 21585.492
$ws.Range("I73").Value = 31604.666
$ws.Range("J73").Value = 4409.7617
$ws.Range("K73").Value = 31604.666
$ws.Range("L73").Value = 4409.7617
$ws.Range("M73").Value = -30668.666
$ws.Range("N73").Value = -6281.7617
$ws.Range("H80").Value = 2396.818
$ws.Range("I80").Value = 2236.5
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 2236.5
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -1238.5
$ws.Range("N80").Value = -5996
$ws.Range("H83").Value = 2396.818
$ws.Range("I83").Value = 2236.5
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 11182.5
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -6190.5
$ws.Range("N83").Value = -29984
$ws.Range("H102").Value = 885.1
$ws.Range("I102").Value = 827.7895
$ws.Range("J102").Value = 984.0909
$ws.Range("K102").Value = 827.7895
$ws.Range("L102").Value = 984.0909
$ws.Range("M102").Value = 794.2105
$ws.Range("N102").Value = -4228.0909

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 36190
$ws.Range("J74").Value = 36190
$ws.Range("L74").Value = 36190
$ws.Range("N74").Value = -38186
$ws.Range("H77").Value = 36190
$ws.Range("J77").Value = 36190
$ws.Range("L77").Value = 108570
$ws.Range("N77").Value = -118554
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents() | Out-Null
$ws.Range("H81").Value = 40000
$ws.Range("J81").Value = 40000
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -41996
$ws.Range("H82").Value = 2489
$ws.Range("I82").Value = 1706.1052
$ws.Range("J82").Value = 4141.778
$ws.Range("K82").Value = 1706.1052
$ws.Range("L82").Value = 4141.778
$ws.Range("M82").Value = -1345.1052
$ws.Range("N82").Value = -4863.778
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents() | Out-Null
$ws.Range("H84").Value = 40000
$ws.Range("J84").Value = 40000
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -129984
$ws.Range("H85").Value = 2489
$ws.Range("I85").Value = 1706.1052
$ws.Range("J85").Value = 4141.778
$ws.Range("K85").Value = 1706.1052
$ws.Range("L85").Value = 4141.778
$ws.Range("M85").Value = -458.1052
$ws.Range("N85").Value = -6637.778

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 18665
$ws.Range("J70").Value = 18665
$ws.Range("L70").Value = 18665
$ws.Range("N70").Value = -19295
$ws.Range("H73").Value = 18665
$ws.Range("J73").Value = 18665
$ws.Range("L73").Value = 18665
$ws.Range("N73").Value = -20849
